# Apply the "rating_summary_template.xlsx" header/column-width update:
#  - Clarify the Begin/End column headers with an explicit timezone label.
#  - Widen the now-longer "Begin"/"End" columns so the new text fits.
#  - Reposition the saved cell selection (cosmetic, matches author's edit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: Begin/End -> Begin(Timezone: UTC)/End(Timezone: UTC) ---
$ws.Range("A1").Value = "Begin(Timezone: UTC)"
$ws.Range("B1").Value = "End(Timezone: UTC)"

# --- Column widths for A and B grew to fit the longer header text ---
# (Excel's ColumnWidth is quantized to whole pixels on the Normal-font grid,
# so we dial in the character width that lands on the requested value.)
$ws.Columns.Item(1).ColumnWidth = 24.428571428571427
$ws.Columns.Item(2).ColumnWidth = 24.857142857142858

# --- Saved selection moved to C15 ---
$ws.Range("C15").Select()
